# Bug fix regarding new permission checks in OpenJDK source.
#
# Sheet "private": rows 78-94 describe individual fields/permission checks.
# A new "java.net.URLConnection / defaultAllowUserInteraction" entry is
# classified as "Irrelevant" (e.g. never used server-side) and moved to the
# top of this block (row 78), the rest of the block shifts down by one, and
# two rows that used to be tagged "Requires change in source code" are
# recategorised: "java.util.Locale / defaultLocale" now only "Requires
# permission" (Neutral) while "java.net.URLConnection / defaultUseCaches"
# keeps needing a source change (Bad/red).
#
# Sheet "reasons": a new reason "Irrelevant" is added with its explanation
# ("E.g. never used on server-side") in column B.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "private"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("private")

$rows = @(
    @{ Row=78; A='java.net.URLConnection';     B='defaultAllowUserInteraction';     C='boolean';                                   D='Irrelevant';                          Style='Normal'  },
    @{ Row=79; A='java.lang.Class';             B='allPermDomain';                   C='java.security.ProtectionDomain';           D='Modification requires permission';    Style='Neutral' },
    @{ Row=80; A='java.lang.SecurityManager';   B='packageDefinitionValid';          C='boolean';                                   D='Modification requires permission';    Style='Neutral' },
    @{ Row=81; A='java.lang.SecurityManager';   B='packageAccessValid';              C='boolean';                                   D='Modification requires permission';    Style='Neutral' },
    @{ Row=82; A='java.lang.System';            B='props';                           C='java.util.Properties';                     D='Modification requires permission';    Style='Neutral' },
    @{ Row=83; A='java.lang.Thread';            B='defaultUncaughtExceptionHandler'; C='java.lang.Thread$UncaughtExceptionHandler'; D='Modification requires permission';   Style='Neutral' },
    @{ Row=84; A='java.net.InetAddress';        B='cachedLocalHost';                 C='java.net.InetAddress$CachedLocalHost';     D='Modification requires permission';    Style='Neutral' },
    @{ Row=85; A='java.net.ServerSocket';       B='factory';                         C='java.net.SocketImplFactory';               D='Modification requires permission';    Style='Neutral' },
    @{ Row=86; A='java.net.Socket';             B='factory';                         C='java.net.SocketImplFactory';               D='Modification requires permission';    Style='Neutral' },
    @{ Row=87; A='java.net.URL';                B='factory';                         C='java.net.URLStreamHandlerFactory';         D='Modification requires permission';    Style='Neutral' },
    @{ Row=88; A='java.net.URLConnection';      B='factory';                         C='java.net.ContentHandlerFactory';           D='Modification requires permission';    Style='Neutral' },
    @{ Row=89; A='java.net.URLConnection';      B='fileNameMap';                     C='java.net.FileNameMap';                     D='Modification requires permission';    Style='Neutral' },
    @{ Row=90; A='java.util.Locale';            B='defaultDisplayLocale';            C='java.util.Locale';                         D='Modification requires permission';    Style='Neutral' },
    @{ Row=91; A='java.util.Locale';            B='defaultFormatLocale';             C='java.util.Locale';                         D='Modification requires permission';    Style='Neutral' },
    @{ Row=92; A='java.util.Locale';            B='defaultLocale';                   C='java.util.Locale';                         D='Modification requires permission';    Style='Neutral' },
    @{ Row=93; A='java.util.TimeZone';          B='defaultTimeZone';                 C='java.util.TimeZone';                       D='Modification requires permission';    Style='Neutral' },
    @{ Row=94; A='java.net.URLConnection';      B='defaultUseCaches';                C='boolean';                                   D='Requires change in source code';      Style='Bad'     }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 4).Style = $r.Style
}

$ws.Activate()
[void]$ws.Range("D94").Select()

# ---------------------------------------------------------------------
# Sheet 5: "reasons" - register the new "Irrelevant" reason + explanation
# ---------------------------------------------------------------------
$reasons = $wb.Worksheets.Item("reasons")
$reasons.Range("A10").Value = 'Irrelevant'
$reasons.Range("B10").Value = 'E.g. never used on server-side'

$reasons.Activate()
[void]$reasons.Range("B10").Select()

# restore the originally active sheet/tab
$ws.Activate()
